# Adapt column header formatting to respective input file names (#7)
# - rename the "_old" / "_new" header-name suffixes to "_FV2410" / "_FV2504"
# - turn the data range A1:U73 into an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1

# --- rename header cell text: "<Name>_old" -> "<Name>_FV2410" -----------
$oldNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

foreach ($name in $oldNames) {
    [void]$ws.Cells.Replace("$name`_old", "$name`_FV2410", $xlWhole)
}

# --- rename header cell text: "<Name>_new" -> "<Name>_FV2504" -----------
foreach ($name in $oldNames) {
    [void]$ws.Cells.Replace("$name`_new", "$name`_FV2504", $xlWhole)
}

# --- turn A1:U73 into a proper Excel table; column names are picked up
#     straight from the (now renamed) header row
$headerRange = $ws.Range("A1:U73")
$tbl = $ws.ListObjects.Add(1, $headerRange, $null, 1)
$tbl.Name = "Table1"

# --- freeze the header row (split below row 1, top-left of scrolling pane A2)
$ws.Activate()
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
